$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns before column D (old D:K shift to F:M)
$ws.Range("D:E").Insert()

# Copy number formats/styles from column F (the old column D, now shifted) onto
# the two new columns so the new D/E cells inherit the correct per-row style
# (date style for header rows, numeric style for data rows) instead of the
# default style picked up from column C during the insert. Bound the range to
# the sheet's used rows so the sheet's dimension/used-range isn't blown out to
# the full 1,048,576-row column.
$ws.Range("F5:F102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Approximate the original bestFit column widths for the two new columns
$ws.Range("D5:E102").ColumnWidth = 14.7109375

# Populate the new D/E columns (most recent two quarters) with the reported
# financial figures for each labeled row.
$data = @(
  @(7, 43465, 43373),
  @(8, 20900, 20300),
  @(9, 15200, 14300),
  @(10, 5700, 6000),
  @(12, 1100, 1100),
  @(13, 0, 0),
  @(14, "NA", "NA"),
  @(15, 100, 100),
  @(17, 19900, 18800),
  @(18, 1000, 1500),
  @(20, 100, 100),
  @(21, 1700, 2100),
  @(22, "NA", 0),
  @(23, 1200, 1500),
  @(24, -18600, 100),
  @(25, 0, 0),
  @(26, 19800, 1400),
  @(27, 19700, 1400),
  @(28, 0, 0),
  @(29, "NA", "NA"),
  @(30, 0, 0),
  @(31, 0, 0),
  @(32, -100, -100),
  @(33, 19700, 1400),
  @(34, 0, 0),
  @(35, 19700, 1400),
  @(38, 43465, 43373),
  @(41, 25600, 25000),
  @(42, 0, 0),
  @(43, 16000, 14500),
  @(44, 22800, 23100),
  @(45, 2800, 3300),
  @(46, 67200, 66000),
  @(47, 0, 0),
  @(48, 10700, 8800),
  @(49, 26600, 26900),
  @(50, 0, 0),
  @(51, 0, 0),
  @(52, 15500, 100),
  @(53, 0, 0),
  @(54, 120100, 101800),
  @(57, 9900, 7300),
  @(58, 0, 0),
  @(59, 5000, 5100),
  @(60, 14900, 12400),
  @(61, 0, 0),
  @(62, 600, 3900),
  @(63, 0, 0),
  @(64, 0, 0),
  @(65, 0, 0),
  @(66, 15500, 16200),
  @(68, 0, 0),
  @(69, 0, 0),
  @(70, 0, 0),
  @(71, 0, 0),
  @(72, -58000, -77700),
  @(73, 0, 0),
  @(74, 0, 0),
  @(75, 0, 0),
  @(76, 104600, 85600),
  @(77, 0, 0),
  @(80, 43465, 43373),
  @(81, 19700, 1400),
  @(83, 600, 600),
  @(84, 0, 0),
  @(85, 0, 0),
  @(86, 0, 0),
  @(87, 0, 0),
  @(88, 0, 0),
  @(89, 3400, 5700),
  @(91, -2200, -1000),
  @(92, 0, 0),
  @(93, 0, 0),
  @(94, -2200, -1000),
  @(96, 0, 0),
  @(97, 0, 0),
  @(98, 0, 0),
  @(99, 0, 0),
  @(100, -600, 100),
  @(101, -100, -100),
  @(102, 500, 4700),
)

foreach ($item in $data) {
  $r = $item[0]
  $ws.Cells.Item($r, 4).Value = $item[1]
  $ws.Cells.Item($r, 5).Value = $item[2]
}
